# Auto-generated Excel COM-interop script
# Applies the "Update countries & provincias Spain" data refresh:
#  - updates the "Datos actualizados..." timestamp in A1
#  - refreshes Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes for the countries whose counts changed
#  - re-sorts the 3 country pairs that swapped rank (by Casos totales, column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Julio de 2020 a las 17:53"

# --- Re-rank the country rows whose order changed ---
# Mozambique overtook Suazilandia (row 136/137)
$ws.Cells.Item(136, 1).Value = "Mozambique"
$ws.Cells.Item(137, 1).Value = "Suazilandia"

# Montenegro overtook Liberia (row 143/144)
$ws.Cells.Item(143, 1).Value = "Montenegro"
$ws.Cells.Item(144, 1).Value = "Liberia"

# Islas Malvinas swapped with Groenlandia (row 209/210, tied totals)
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(210, 1).Value = "Groenlandia"

# --- Refresh the numeric statistics (Casos totales .. Muertes) ---
# Row 4
$ws.Cells.Item(4, 2).Value = 3114576
$ws.Cells.Item(4, 3).Value = 17492
$ws.Cells.Item(4, 4).Value = 1355855
$ws.Cells.Item(4, 5).Value = 1624580
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 169
$ws.Cells.Item(4, 8).Value = 134141

# Row 6
$ws.Cells.Item(6, 2).Value = 760761
$ws.Cells.Item(6, 3).Value = 17280
$ws.Cells.Item(6, 4).Value = 470144
$ws.Cells.Item(6, 5).Value = 269545
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 419
$ws.Cells.Item(6, 8).Value = 21072

# Row 10
$ws.Cells.Item(10, 2).Value = 299593
$ws.Cells.Item(10, 3).Value = 383
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 4
$ws.Cells.Item(10, 8).Value = 28396

# Row 11
$ws.Cells.Item(11, 2).Value = 286979
$ws.Cells.Item(11, 3).Value = 630
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 126
$ws.Cells.Item(11, 8).Value = 44517

# Row 19
$ws.Cells.Item(19, 2).Value = 198561
$ws.Cells.Item(19, 3).Value = 206
$ws.Cells.Item(19, 4).Value = 182700
$ws.Cells.Item(19, 5).Value = 6756
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 2
$ws.Cells.Item(19, 8).Value = 9105

# Row 41
$ws.Cells.Item(41, 2).Value = 45298
$ws.Cells.Item(41, 3).Value = 158
$ws.Cells.Item(41, 4).Value = 41323
$ws.Cells.Item(41, 5).Value = 3949
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 26

# Row 42
$ws.Cells.Item(42, 2).Value = 44859
$ws.Cells.Item(42, 3).Value = 443
$ws.Cells.Item(42, 4).Value = 29714
$ws.Cells.Item(42, 5).Value = 13514
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = 1631

# Row 45
$ws.Cells.Item(45, 2).Value = 39588
$ws.Cells.Item(45, 3).Value = 1158
$ws.Cells.Item(45, 4).Value = 20056
$ws.Cells.Item(45, 5).Value = 18703
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 8
$ws.Cells.Item(45, 8).Value = 829

# Row 58
$ws.Cells.Item(58, 2).Value = 21916
$ws.Cells.Item(58, 3).Value = 542
$ws.Cells.Item(58, 4).Value = 13100
$ws.Cells.Item(58, 5).Value = 8542
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 9
$ws.Cells.Item(58, 8).Value = 274

# Row 61
$ws.Cells.Item(61, 2).Value = 18471
$ws.Cells.Item(61, 3).Value = 330
$ws.Cells.Item(61, 4).Value = 11549
$ws.Cells.Item(61, 5).Value = 6308
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 11
$ws.Cells.Item(61, 8).Value = 614

# Row 69
$ws.Cells.Item(69, 2).Value = 12775
$ws.Cells.Item(69, 3).Value = 90
$ws.Cells.Item(69, 4).Value = 8005
$ws.Cells.Item(69, 5).Value = 4419
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 351

# Row 73
$ws.Cells.Item(73, 2).Value = 8950
$ws.Cells.Item(73, 3).Value = 3
$ws.Cells.Item(73, 4).Value = 8138
$ws.Cells.Item(73, 5).Value = 561
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 251

# Row 99
$ws.Cells.Item(99, 2).Value = 3622
$ws.Cells.Item(99, 3).Value = 33
$ws.Cells.Item(99, 4).Value = 1374
$ws.Cells.Item(99, 5).Value = 2055
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 193

# Row 109
$ws.Cells.Item(109, 2).Value = 2501
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 2158
$ws.Cells.Item(109, 5).Value = 330
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = 13

# Row 112
$ws.Cells.Item(112, 2).Value = 2093
$ws.Cells.Item(112, 3).Value = 12
$ws.Cells.Item(112, 4).Value = 1967
$ws.Cells.Item(112, 5).Value = 115
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 11

# Row 115
$ws.Cells.Item(115, 2).Value = 1946
$ws.Cells.Item(115, 3).Value = 39
$ws.Cells.Item(115, 4).Value = 1368
$ws.Cells.Item(115, 5).Value = 542
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 36

# Row 136
$ws.Cells.Item(136, 2).Value = 1071
$ws.Cells.Item(136, 3).Value = 31
$ws.Cells.Item(136, 4).Value = 337
$ws.Cells.Item(136, 5).Value = 726
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 8

# Row 137
$ws.Cells.Item(137, 2).Value = 1056
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 570
$ws.Cells.Item(137, 5).Value = 472
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 14

# Row 139
$ws.Cells.Item(139, 2).Value = 1003
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 861
$ws.Cells.Item(139, 5).Value = 89
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 53

# Row 143
$ws.Cells.Item(143, 2).Value = 960
$ws.Cells.Item(143, 3).Value = 53
$ws.Cells.Item(143, 4).Value = 320
$ws.Cells.Item(143, 5).Value = 623
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 17

# Row 144
$ws.Cells.Item(144, 2).Value = 917
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 394
$ws.Cells.Item(144, 5).Value = 482
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 41
